$wb = $excel.ActiveWorkbook

# --- Step 1: adjust two values in the existing May2020 sheet ---
$may = $wb.Worksheets.Item("May2020")
$may.Range("D25").Value = 0
$may.Range("D38").Value = 2

# --- Step 2: add the new Jun2020 sheet right after May2020 ---
$new = $wb.Worksheets.Add($null, $may)
$new.Name = "Jun2020"

# --- Step 3: populate header row (bold, bordered, centered/top) ---
$new.Cells.Item(1, 1).Value = "name"
$new.Cells.Item(1, 2).Value = "team"
$new.Cells.Item(1, 3).Value = "IOMP-MT"
$new.Cells.Item(1, 4).Value = "IOMP-CT"
$headerRange = $new.Range("A1:D1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108
$headerRange.VerticalAlignment = -4160
$headerRange.Borders.LineStyle = 1

# --- Step 4: populate data rows 2-54 ---
$new.Cells.Item(2, 1).Value = "Amy"
$new.Cells.Item(2, 2).Value = "admin"
$new.Cells.Item(2, 3).Value = 0
$new.Cells.Item(2, 4).Value = 1
$new.Cells.Item(3, 1).Value = "Anj"
$new.Cells.Item(3, 2).Value = "MT"
$new.Cells.Item(3, 3).Value = 1
$new.Cells.Item(3, 4).Value = 1
$new.Cells.Item(4, 1).Value = "Anne"
$new.Cells.Item(4, 2).Value = "CT"
$new.Cells.Item(4, 3).Value = 1
$new.Cells.Item(4, 4).Value = 1
$new.Cells.Item(5, 1).Value = "Ardeth"
$new.Cells.Item(5, 2).Value = "CT"
$new.Cells.Item(5, 3).Value = 1
$new.Cells.Item(5, 4).Value = 1
$new.Cells.Item(6, 1).Value = "Arnel"
$new.Cells.Item(6, 2).Value = "MT"
$new.Cells.Item(6, 3).Value = 1
$new.Cells.Item(6, 4).Value = 1
$new.Cells.Item(7, 1).Value = "Brain"
$new.Cells.Item(7, 2).Value = "MT"
$new.Cells.Item(7, 3).Value = 1
$new.Cells.Item(7, 4).Value = 1
$new.Cells.Item(8, 1).Value = "Camille"
$new.Cells.Item(8, 2).Value = "CT"
$new.Cells.Item(8, 3).Value = 1
$new.Cells.Item(8, 4).Value = 1
$new.Cells.Item(9, 1).Value = "Carla"
$new.Cells.Item(9, 2).Value = "CT"
$new.Cells.Item(9, 3).Value = 1
$new.Cells.Item(9, 4).Value = 1
$new.Cells.Item(10, 1).Value = "Carlo"
$new.Cells.Item(10, 2).Value = "CT"
$new.Cells.Item(10, 3).Value = 1
$new.Cells.Item(10, 4).Value = 1
$new.Cells.Item(11, 1).Value = "Cath"
$new.Cells.Item(11, 2).Value = "MT"
$new.Cells.Item(11, 3).Value = 2
$new.Cells.Item(11, 4).Value = 1
$new.Cells.Item(12, 1).Value = "Chad"
$new.Cells.Item(12, 2).Value = "MT"
$new.Cells.Item(12, 3).Value = 1
$new.Cells.Item(12, 4).Value = 1
$new.Cells.Item(13, 1).Value = "Daisy"
$new.Cells.Item(13, 2).Value = "admin"
$new.Cells.Item(13, 3).Value = 0
$new.Cells.Item(13, 4).Value = 1
$new.Cells.Item(14, 1).Value = "Dan"
$new.Cells.Item(14, 2).Value = "CT"
$new.Cells.Item(14, 3).Value = 1
$new.Cells.Item(14, 4).Value = 2
$new.Cells.Item(15, 1).Value = "David"
$new.Cells.Item(15, 2).Value = "MT"
$new.Cells.Item(15, 3).Value = 1
$new.Cells.Item(15, 4).Value = 1
$new.Cells.Item(16, 1).Value = "Drew"
$new.Cells.Item(16, 2).Value = "MT"
$new.Cells.Item(16, 3).Value = 1
$new.Cells.Item(16, 4).Value = 1
$new.Cells.Item(17, 1).Value = "Eunice"
$new.Cells.Item(17, 2).Value = "CT"
$new.Cells.Item(17, 3).Value = 1
$new.Cells.Item(17, 4).Value = 2
$new.Cells.Item(18, 1).Value = "Gene"
$new.Cells.Item(18, 2).Value = "MT"
$new.Cells.Item(18, 3).Value = 2
$new.Cells.Item(18, 4).Value = 1
$new.Cells.Item(19, 1).Value = "Harle"
$new.Cells.Item(19, 2).Value = "MT"
$new.Cells.Item(19, 3).Value = 2
$new.Cells.Item(19, 4).Value = 1
$new.Cells.Item(20, 1).Value = "Harry"
$new.Cells.Item(20, 2).Value = "CT"
$new.Cells.Item(20, 3).Value = 1
$new.Cells.Item(20, 4).Value = 2
$new.Cells.Item(21, 1).Value = "Issa"
$new.Cells.Item(21, 2).Value = "CT"
$new.Cells.Item(21, 3).Value = 1
$new.Cells.Item(21, 4).Value = 1
$new.Cells.Item(22, 1).Value = "JK"
$new.Cells.Item(22, 2).Value = "MT"
$new.Cells.Item(22, 3).Value = 1
$new.Cells.Item(22, 4).Value = 1
$new.Cells.Item(23, 1).Value = "James"
$new.Cells.Item(23, 2).Value = "MT"
$new.Cells.Item(23, 3).Value = 2
$new.Cells.Item(23, 4).Value = 1
$new.Cells.Item(24, 1).Value = "Janine"
$new.Cells.Item(24, 2).Value = "MT"
$new.Cells.Item(24, 3).Value = 1
$new.Cells.Item(24, 4).Value = 1
$new.Cells.Item(25, 1).Value = "Jec"
$new.Cells.Item(25, 2).Value = "CT"
$new.Cells.Item(25, 3).Value = 1
$new.Cells.Item(25, 4).Value = 2
$new.Cells.Item(26, 1).Value = "Jes"
$new.Cells.Item(26, 2).Value = "CT"
$new.Cells.Item(26, 3).Value = 2
$new.Cells.Item(26, 4).Value = 1
$new.Cells.Item(27, 1).Value = "Jhoanne"
$new.Cells.Item(27, 2).Value = "admin"
$new.Cells.Item(27, 3).Value = 0
$new.Cells.Item(27, 4).Value = 1
$new.Cells.Item(28, 1).Value = "Johann"
$new.Cells.Item(28, 2).Value = "MT"
$new.Cells.Item(28, 3).Value = 1
$new.Cells.Item(28, 4).Value = 1
$new.Cells.Item(29, 1).Value = "John"
$new.Cells.Item(29, 2).Value = "MT"
$new.Cells.Item(29, 3).Value = 2
$new.Cells.Item(29, 4).Value = 1
$new.Cells.Item(30, 1).Value = "Julius"
$new.Cells.Item(30, 2).Value = "CT"
$new.Cells.Item(30, 3).Value = 1
$new.Cells.Item(30, 4).Value = 1
$new.Cells.Item(31, 1).Value = "Kate"
$new.Cells.Item(31, 2).Value = "MT"
$new.Cells.Item(31, 3).Value = 1
$new.Cells.Item(31, 4).Value = 1
$new.Cells.Item(32, 1).Value = "Ken"
$new.Cells.Item(32, 2).Value = "CT"
$new.Cells.Item(32, 3).Value = 1
$new.Cells.Item(32, 4).Value = 1
$new.Cells.Item(33, 1).Value = "Kennex"
$new.Cells.Item(33, 2).Value = "MT"
$new.Cells.Item(33, 3).Value = 1
$new.Cells.Item(33, 4).Value = 1
$new.Cells.Item(34, 1).Value = "Kevin"
$new.Cells.Item(34, 2).Value = "MT"
$new.Cells.Item(34, 3).Value = 2
$new.Cells.Item(34, 4).Value = 1
$new.Cells.Item(35, 1).Value = "Lem"
$new.Cells.Item(35, 2).Value = "CT"
$new.Cells.Item(35, 3).Value = 1
$new.Cells.Item(35, 4).Value = 1
$new.Cells.Item(36, 1).Value = "Louie"
$new.Cells.Item(36, 2).Value = "MT"
$new.Cells.Item(36, 3).Value = 1
$new.Cells.Item(36, 4).Value = 1
$new.Cells.Item(37, 1).Value = "Luz"
$new.Cells.Item(37, 2).Value = "MT"
$new.Cells.Item(37, 3).Value = 2
$new.Cells.Item(37, 4).Value = 1
$new.Cells.Item(38, 1).Value = "Meryll"
$new.Cells.Item(38, 2).Value = "MT"
$new.Cells.Item(38, 3).Value = 1
$new.Cells.Item(38, 4).Value = 1
$new.Cells.Item(39, 1).Value = "Momay"
$new.Cells.Item(39, 2).Value = "CT"
$new.Cells.Item(39, 3).Value = 1
$new.Cells.Item(39, 4).Value = 1
$new.Cells.Item(40, 1).Value = "Morgan"
$new.Cells.Item(40, 2).Value = "MT"
$new.Cells.Item(40, 3).Value = 1
$new.Cells.Item(40, 4).Value = 1
$new.Cells.Item(41, 1).Value = "Nathan"
$new.Cells.Item(41, 2).Value = "MT"
$new.Cells.Item(41, 3).Value = 2
$new.Cells.Item(41, 4).Value = 1
$new.Cells.Item(42, 1).Value = "Nichole"
$new.Cells.Item(42, 2).Value = "MT"
$new.Cells.Item(42, 3).Value = 2
$new.Cells.Item(42, 4).Value = 1
$new.Cells.Item(43, 1).Value = "Nora"
$new.Cells.Item(43, 2).Value = "CT"
$new.Cells.Item(43, 3).Value = 1
$new.Cells.Item(43, 4).Value = 2
$new.Cells.Item(44, 1).Value = "Oscar"
$new.Cells.Item(44, 2).Value = "CT"
$new.Cells.Item(44, 3).Value = 1
$new.Cells.Item(44, 4).Value = 1
$new.Cells.Item(45, 1).Value = "Pau"
$new.Cells.Item(45, 2).Value = "CT"
$new.Cells.Item(45, 3).Value = 1
$new.Cells.Item(45, 4).Value = 1
$new.Cells.Item(46, 1).Value = "Reyn"
$new.Cells.Item(46, 2).Value = "CT"
$new.Cells.Item(46, 3).Value = 1
$new.Cells.Item(46, 4).Value = 2
$new.Cells.Item(47, 1).Value = "Rodney"
$new.Cells.Item(47, 2).Value = "CT"
$new.Cells.Item(47, 3).Value = 1
$new.Cells.Item(47, 4).Value = 1
$new.Cells.Item(48, 1).Value = "Roy"
$new.Cells.Item(48, 2).Value = "MT"
$new.Cells.Item(48, 3).Value = 1
$new.Cells.Item(48, 4).Value = 1
$new.Cells.Item(49, 1).Value = "Sky"
$new.Cells.Item(49, 2).Value = "MT"
$new.Cells.Item(49, 3).Value = 2
$new.Cells.Item(49, 4).Value = 1
$new.Cells.Item(50, 1).Value = "Tina"
$new.Cells.Item(50, 2).Value = "MT"
$new.Cells.Item(50, 3).Value = 1
$new.Cells.Item(50, 4).Value = 1
$new.Cells.Item(51, 1).Value = "Tine"
$new.Cells.Item(51, 2).Value = "CT"
$new.Cells.Item(51, 3).Value = 1
$new.Cells.Item(51, 4).Value = 1
$new.Cells.Item(52, 1).Value = "Tintin"
$new.Cells.Item(52, 2).Value = "admin"
$new.Cells.Item(52, 3).Value = 0
$new.Cells.Item(52, 4).Value = 1
$new.Cells.Item(53, 1).Value = "Troy"
$new.Cells.Item(53, 2).Value = "MT"
$new.Cells.Item(53, 3).Value = 1
$new.Cells.Item(53, 4).Value = 1
$new.Cells.Item(54, 1).Value = "Web"
$new.Cells.Item(54, 2).Value = "CT"
$new.Cells.Item(54, 3).Value = 1
$new.Cells.Item(54, 4).Value = 2

Write-Output "done"
